# ex9.1.2(Linear) Strong Stationary generator alpha non zero
# "volver a generar problemas cuadraticos y lineales"
# Re-generated values for the follower restrictions, the modified point,
# vec_bf, vec_BF and vec_alpha sheets.

$wb = $excel.ActiveWorkbook

# ---- Restricciones_del_follower ------------------------------------------------
$wsFollower = $wb.Worksheets.Item(3)   # "Restricciones_del_follower"

# Row 2 (J_0_L0_v)
$wsFollower.Range("A2").Value = "-2.536821705426357 + 0.4263565891472869y"
$wsFollower.Range("B2").NumberFormat = "@"
$wsFollower.Range("B2").Value = "2.536821705426357"
$wsFollower.Range("D2").NumberFormat = "@"
$wsFollower.Range("D2").Value = "0.37"
$wsFollower.Range("E2").NumberFormat = "@"
$wsFollower.Range("E2").Value = "6.5"
$wsFollower.Range("F2").NumberFormat = "@"
$wsFollower.Range("F2").Value = "1.1"

# Row 3 (J_0_L0_v)
$wsFollower.Range("A3").Value = "-12.868992248062016 - x + 3.062015503875969y"
$wsFollower.Range("B3").NumberFormat = "@"
$wsFollower.Range("B3").Value = "9.868992248062016"
$wsFollower.Range("D3").NumberFormat = "@"
$wsFollower.Range("D3").Value = "0.44"
$wsFollower.Range("E3").NumberFormat = "@"
$wsFollower.Range("E3").Value = "7.7"
$wsFollower.Range("F3").NumberFormat = "@"
$wsFollower.Range("F3").Value = "7.9"

# Row 4 (J_0_LP_v)
$wsFollower.Range("A4").Value = "13.452000000000005 + x - 3.16y"
$wsFollower.Range("B4").NumberFormat = "@"
$wsFollower.Range("B4").Value = "-25.452000000000005"
$wsFollower.Range("D4").NumberFormat = "@"
$wsFollower.Range("D4").Value = "0.0"
$wsFollower.Range("E4").NumberFormat = "@"
$wsFollower.Range("E4").Value = "1.7000000000000002"
$wsFollower.Range("F4").NumberFormat = "@"
$wsFollower.Range("F4").Value = "5.5"

# Row 5 (J_Ne_L0_v)
$wsFollower.Range("A5").Value = "-34.392248062015504 + 4x + 2.0155038759689923y"
$wsFollower.Range("B5").NumberFormat = "@"
$wsFollower.Range("B5").Value = "21.392248062015504"
$wsFollower.Range("D5").NumberFormat = "@"
$wsFollower.Range("D5").Value = "0.07"
$wsFollower.Range("E5").NumberFormat = "@"
$wsFollower.Range("E5").Value = "5.699999999999999"
$wsFollower.Range("F5").NumberFormat = "@"
$wsFollower.Range("F5").Value = "5.2"

# ---- Punto_modificado -----------------------------------------------------
$wsPunto = $wb.Worksheets.Item(4)   # "Punto_modificado"
$wsPunto.Range("A2").NumberFormat = "@"
$wsPunto.Range("A2").Value = "5.35"
$wsPunto.Range("B2").NumberFormat = "@"
$wsPunto.Range("B2").Value = "5.95"

# ---- Vector_bf / Vector_BF ---------------------------------------------------
# NOTE: sheet names only differ by case ("Vector_bf" vs "Vector_BF"), and
# Worksheets.Item(<name>) resolves case-insensitively (both names would hit
# the same sheet) - so address these two by their 1-based tab position
# instead of by name to land on the correct sheet.
$wsBf = $wb.Worksheets.Item(5)   # "Vector_bf"
$wsBf.Range("A2").NumberFormat = "@"
$wsBf.Range("A2").Value = "-2.646124031007752"

$wsBF = $wb.Worksheets.Item(6)   # "Vector_BF"
$wsBF.Range("A2").NumberFormat = "@"
$wsBF.Range("A2").Value = "-15.799999999999997"
$wsBF.Range("A3").NumberFormat = "@"
$wsBF.Range("A3").Value = "-29.465209302325576"

# ---- Vector_Alpha -----------------------------------------------------------
# A2 here is a genuine numeric cell (not text), unlike all the cells above.
$wsAlpha = $wb.Worksheets.Item(7) # "Vector_Alpha"
$wsAlpha.Range("A2").Value = 2.58

Write-Output "edits applied"
